$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Updated "VENTA" (sales) figures
$ws.Range("D2").Value = 73996.77
$ws.Range("D3").Value = 3361.24

# Derived "POR CUMPLIR" values (PRESUPUESTO - VENTA)
$ws.Range("E2").Value = -73996.77
$ws.Range("E3").Value = 12110.3193

# Derived "CUMPLIMIENTO" percentages (VENTA / PRESUPUESTO)
$ws.Range("F3").Value = 0.2172528272570432

# Totals row (row 4)
$ws.Range("D4").Value = 77358.01000000001
$ws.Range("E4").Value = -61886.4507
$ws.Range("F4").Value = 5.000013799514054
